$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new values
$ws.Range("B2").Value = "5002019-61.2019.8.21.0067"
$ws.Range("C2").Value = "0064471-68.2019.8.21.9000"
$ws.Range("D2").Value = "CIV.04109.02"
$ws.Range("E2").Value = "originario_principal"

# Delete rows 3 through 5 (entire rows), shifting cells up
$ws.Range("A3:E5").Delete()
